$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.924.95"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "1.860.42"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  -1.74%  "

$ws.Range("D5").Value = "'320.77"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("D7").Value = "'0.4369"
$ws.Range("E7").Value = "  -1.13%  "

$ws.Range("D8").Value = "'0.3785"
$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("D9").Value = "'0.07439"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("D10").Value = "'0.8860"
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("D11").Value = "'21.70"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "1.859.40"
$ws.Range("E12").Value = "  -0.59%  "

$ws.Range("D13").Value = "'6.795"
$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").Value = "'5.503"
$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("D15").Value = "'0.07149"
$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").Value = "'88.64"
$ws.Range("E16").Value = "  +5.68%  "

$ws.Range("D17").Value = "'1.021"
$ws.Range("E17").Value = "  -1.64%  "

$ws.Range("D18").Value = "'0.000009038"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").Value = "'1.017"
$ws.Range("E19").Value = "  -1.53%  "

$ws.Range("D20").Value = "'15.50"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").Value = "27.891.55"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").Value = "'5.272"
$ws.Range("E22").Value = "  -0.72%  "

$ws.Range("D23").Value = "'11.19"
$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("D24").Value = "2.080.19"
$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("E25").Value = "  +4.88%  "

$ws.Range("D26").Value = "'157.03"
$ws.Range("E26").Value = "  -0.94%  "

$ws.Range("D27").Value = "'18.73"
$ws.Range("E27").Value = "  -0.68%  "

$ws.Range("D28").Value = "'5.450"
$ws.Range("E28").Value = "  +2.02%  "

$ws.Range("D29").Value = "'2.004"
$ws.Range("E29").Value = "  +0.50%  "

$ws.Range("D30").Value = "'120.64"
$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("D31").Value = "'0.09037"
$ws.Range("E31").Value = "  -0.67%  "

$ws.Range("D32").Value = "'1.237"
$ws.Range("E32").Value = "  +1.78%  "

$ws.Range("D33").Value = "'0.7739"
$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("D34").Value = "'4.580"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").Value = "'2.987"
$ws.Range("E35").Value = "  +2.48%  "

$ws.Range("E36").Value = "  -1.54%  "

$ws.Range("D37").Value = "'1.141"
$ws.Range("E37").Value = "  -1.58%  "

$ws.Range("D38").Value = "'0.01981"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("D39").Value = "'0.05318"
$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").Value = "'2.879"
$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("D41").Value = "'0.5208"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").Value = "'6.992"
$ws.Range("E42").Value = "  +2.08%  "

$ws.Range("D43").Value = "'0.1680"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("D44").Value = "'8.758"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "'110.50"
$ws.Range("E45").Value = "  +0.77%  "

$ws.Range("D46").Value = "'10.75"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.715"
$ws.Range("E47").Value = "  -1.10%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4749"
$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("D49").Value = "'1.018"
$ws.Range("E49").Value = "  -1.65%  "

$ws.Range("D50").Value = "'0.06478"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("D51").Value = "'1.855"
$ws.Range("E51").Value = "  -0.70%  "
